# Restore revision: update "R30" rule's "Integer min" value from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
